$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) Insert two additional data rows into the workers table so it
#    grows from 6 rows (16-21) to 8 rows (16-23). We insert before
#    the former last row (21) so that the last row (with its
#    special "closing" bottom-border styling) ends up at row 23.
# ---------------------------------------------------------------
$ws.Rows("21:22").Insert()

# Copy the formatting (styles / number formats / borders) of row 20
# (a normal, non-closing table row) into the two newly inserted rows.
$ws.Range("B20:J20").Copy() | Out-Null
$ws.Range("B21:J21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B20:J20").Copy() | Out-Null
$ws.Range("B22:J22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 2) Rewrite the workers table (rows 16-23) with the new data set:
#    previous accounts for IVAN MENDOZA RAMIREZ, LUIS CARLOS
#    GUTIERREZ VEGA and JHONATHAN RECUERO MORELO are kept, a new
#    worker (JOSE VICENTE GUZMAN) is added, and each worker now has
#    two rows (periods 1610 and 1609, in that order).
# ---------------------------------------------------------------
$workers = @(
    @{ Doc = "1129495372"; Name = "IVAN MENDOZA RAMIREZ" },
    @{ Doc = "1045309580"; Name = "JOSE VICENTE GUZMAN" },
    @{ Doc = "73559861";   Name = "LUIS CARLOS GUTIERREZ VEGA" },
    @{ Doc = "73212463";   Name = "JHONATHAN RECUERO MORELO" }
)
$periods = @("1610", "1609")

$row = 16
foreach ($worker in $workers) {
    foreach ($period in $periods) {
        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $worker.Doc
        $ws.Cells.Item($row, 4).Value = $worker.Name
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = 27578
        $ws.Cells.Item($row, 7).Value = 689455
        $row = $row + 1
    }
}

# ---------------------------------------------------------------
# 3) Update the summary figures above the table.
# ---------------------------------------------------------------
# Total "VALOR MORA" (sum of the Valor Mora column for the 8 rows)
$ws.Range("E11").Value = 220624
# "Cant. Trabajadores" now lists 4 distinct workers
$ws.Range("C13").Value = 4

$wb.Save()
